# Refresh the cryptos table with the latest scraped snapshot (coinranking.com).
# Most rows only get fresh Price/Volume(1h) figures; a handful of rows were
# re-ranked by the scraper, so their Coin/Link/Price/Volume move to a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Some "prices" are plain decimal-looking strings (e.g. "1.00", "0.0420")
    # that must stay text so formatting like trailing zeros survives - a bare
    # assignment would let Excel reinterpret them as numbers. A leading
    # apostrophe forces text entry; reset the style afterwards so the cell
    # doesn't keep the quote-prefix formatting Excel applies automatically.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "96.055.68"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "3.489.68"
$ws.Range("E3").Value = "  +4.94%  "

$ws.Range("E4").Value = "  +0.22%  "

Set-TextValue $ws.Range("D5") "242.12"
$ws.Range("E5").Value = "  -2.34%  "

Set-TextValue $ws.Range("D6") "647.63"
$ws.Range("E6").Value = "  -0.56%  "

Set-TextValue $ws.Range("D7") "1.49"
$ws.Range("E7").Value = "  +9.99%  "

Set-TextValue $ws.Range("D8") "0.414"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("E10").Value = "  +2.18%  "

$ws.Range("D11").Value = "3.486.31"
$ws.Range("E11").Value = "  +4.91%  "

Set-TextValue $ws.Range("D12") "42.78"
$ws.Range("E12").Value = "  +7.97%  "

Set-TextValue $ws.Range("D13") "0.199"
$ws.Range("E13").Value = "  -2.48%  "

Set-TextValue $ws.Range("D14") "6.16"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").Value = "96.054.69"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").Value = "4.148.36"
$ws.Range("E16").Value = "  +5.32%  "

Set-TextValue $ws.Range("D17") "0.0000254"
$ws.Range("E17").Value = "  +1.75%  "

Set-TextValue $ws.Range("D18") "8.49"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "3.514.81"
$ws.Range("E19").Value = "  +5.90%  "

Set-TextValue $ws.Range("D20") "18.38"
$ws.Range("E20").Value = "  +9.94%  "

Set-TextValue $ws.Range("D21") "11.83"
$ws.Range("E21").Value = "  +13.88%  "

Set-TextValue $ws.Range("D22") "0.525"
$ws.Range("E22").Value = "  +10.31%  "

Set-TextValue $ws.Range("D23") "512.46"
$ws.Range("E23").Value = "  +2.73%  "

Set-TextValue $ws.Range("D24") "3.25"
$ws.Range("E24").Value = "  -2.28%  "

Set-TextValue $ws.Range("D25") "0.0000194"
$ws.Range("E25").Value = "  -1.31%  "

Set-TextValue $ws.Range("D26") "6.63"
$ws.Range("E26").Value = "  +3.11%  "

Set-TextValue $ws.Range("D27") "92.72"
$ws.Range("E27").Value = "  -1.57%  "

Set-TextValue $ws.Range("D28") "12.40"
$ws.Range("E28").Value = "  +4.25%  "

$ws.Range("D29").Value = "3.691.83"
$ws.Range("E29").Value = "  +5.66%  "

Set-TextValue $ws.Range("D30") "11.90"
$ws.Range("E30").Value = "  +10.94%  "

# Row 31 was re-ranked by the scraper; now "Dai"
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D31") "0.995"
$ws.Range("E31").Value = "  -0.56%  "

# Row 32 was re-ranked by the scraper; now "PancakeSwap"
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D32") "2.78"
$ws.Range("E32").Value = "  +12.96%  "

Set-TextValue $ws.Range("D33") "0.139"
$ws.Range("E33").Value = "  -1.75%  "

Set-TextValue $ws.Range("D34") "0.185"
$ws.Range("E34").Value = "  +0.01%  "

Set-TextValue $ws.Range("D35") "31.04"
$ws.Range("E35").Value = "  +11.36%  "

Set-TextValue $ws.Range("D36") "0.578"
$ws.Range("E36").Value = "  +7.03%  "

Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  -0.34%  "

Set-TextValue $ws.Range("D38") "7.84"
$ws.Range("E38").Value = "  +4.88%  "

Set-TextValue $ws.Range("D39") "1.47"
$ws.Range("E39").Value = "  -2.02%  "

Set-TextValue $ws.Range("D40") "520.71"
$ws.Range("E40").Value = "  +3.53%  "

Set-TextValue $ws.Range("D41") "0.152"
$ws.Range("E41").Value = "  +1.22%  "

# Row 42 was re-ranked by the scraper; now "USDe"
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D42") "1.00"
$ws.Range("E42").Value = "  +0.05%  "

# Row 43 was re-ranked by the scraper; now "ARBITRUM"
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D43") "0.918"
$ws.Range("E43").Value = "  +11.79%  "

Set-TextValue $ws.Range("D44") "24.13"
$ws.Range("E44").Value = "  -1.40%  "

Set-TextValue $ws.Range("D45") "1.72"
$ws.Range("E45").Value = "  +7.00%  "

Set-TextValue $ws.Range("D46") "0.0420"
$ws.Range("E46").Value = "  +4.02%  "

# Row 47 was re-ranked by the scraper; now "Filecoin"
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D47") "5.60"
$ws.Range("E47").Value = "  +3.80%  "

# Row 48 was re-ranked by the scraper; now "MantraDAO"
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D48") "3.61"
$ws.Range("E48").Value = "  -1.14%  "

# Row 49 was re-ranked by the scraper; now "Stacks"
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D49") "2.18"
$ws.Range("E49").Value = "  +12.14%  "

# Row 50 was re-ranked by the scraper; now "dogwifhat"
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D50") "3.23"
$ws.Range("E50").Value = "  +3.57%  "

Set-TextValue $ws.Range("D51") "8.24"
$ws.Range("E51").Value = "  -0.52%  "
